$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '64.306.82'
$ws.Cells.Item(2, 5).Value = '  -0.69%  '
$ws.Cells.Item(3, 4).Value = '3.508.47'
$ws.Cells.Item(3, 5).Value = '  +0.17%  '
$ws.Cells.Item(4, 5).Value = '  -0.07%  '
$ws.Cells.Item(5, 4).Value = "'583.91"
$ws.Cells.Item(5, 5).Value = '  -0.46%  '
$ws.Cells.Item(6, 4).Value = "'134.63"
$ws.Cells.Item(6, 5).Value = '  +1.81%  '
$ws.Cells.Item(7, 4).Value = '3.508.58'
$ws.Cells.Item(7, 5).Value = '  +0.25%  '
$ws.Cells.Item(8, 5).Value = '  -0.03%  '
$ws.Cells.Item(9, 5).Value = '  -0.05%  '
$ws.Cells.Item(10, 5).Value = '  +0.63%  '
$ws.Cells.Item(11, 4).Value = "'7.09"
$ws.Cells.Item(11, 5).Value = '  -0.36%  '
$ws.Cells.Item(12, 4).Value = "'0.374"
$ws.Cells.Item(12, 5).Value = '  -2.75%  '
$ws.Cells.Item(13, 4).Value = '4.104.78'
$ws.Cells.Item(13, 5).Value = '  +0.06%  '
$ws.Cells.Item(14, 5).Value = '  +0.16%  '
$ws.Cells.Item(15, 2).Value = 'TRON'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(15, 4).Value = "'0.118"
$ws.Cells.Item(15, 5).Value = '  +1.06%  '
$ws.Cells.Item(16, 2).Value = 'WrappedEther'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(16, 4).Value = '3.506.66'
$ws.Cells.Item(16, 5).Value = '  +0.00%  '
$ws.Cells.Item(17, 2).Value = 'Avalanche'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(17, 4).Value = "'26.57"
$ws.Cells.Item(17, 5).Value = '  -3.93%  '
$ws.Cells.Item(18, 4).Value = '64.319.59'
$ws.Cells.Item(18, 5).Value = '  -0.74%  '
$ws.Cells.Item(19, 5).Value = '  -2.13%  '
$ws.Cells.Item(20, 4).Value = "'13.82"
$ws.Cells.Item(20, 5).Value = '  -2.80%  '
$ws.Cells.Item(21, 4).Value = "'5.58"
$ws.Cells.Item(21, 5).Value = '  -1.37%  '
$ws.Cells.Item(22, 4).Value = "'383.73"
$ws.Cells.Item(22, 5).Value = '  -1.92%  '
$ws.Cells.Item(23, 2).Value = 'WrappedeETH'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(23, 4).Value = '3.649.93'
$ws.Cells.Item(23, 5).Value = '  +0.08%  '
$ws.Cells.Item(24, 2).Value = 'Polygon'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(24, 4).Value = "'0.568"
$ws.Cells.Item(24, 5).Value = '  -1.28%  '
$ws.Cells.Item(25, 4).Value = "'73.95"
$ws.Cells.Item(25, 5).Value = '  -0.29%  '
$ws.Cells.Item(26, 5).Value = '  -0.05%  '
$ws.Cells.Item(27, 5).Value = '  -0.05%  '
$ws.Cells.Item(28, 5).Value = '  +4.90%  '
$ws.Cells.Item(29, 4).Value = "'7.62"
$ws.Cells.Item(29, 5).Value = '  +3.01%  '
$ws.Cells.Item(30, 5).Value = '  -1.08%  '
$ws.Cells.Item(31, 5).Value = '  +0.04%  '
$ws.Cells.Item(32, 4).Value = "'8.29"
$ws.Cells.Item(32, 5).Value = '  +1.15%  '
$ws.Cells.Item(33, 4).Value = "'2.21"
$ws.Cells.Item(33, 5).Value = '  -2.02%  '
$ws.Cells.Item(34, 4).Value = '3.526.98'
$ws.Cells.Item(34, 5).Value = '  +0.49%  '
$ws.Cells.Item(35, 5).Value = '  -0.03%  '
$ws.Cells.Item(36, 5).Value = '  +0.56%  '
$ws.Cells.Item(37, 4).Value = "'23.53"
$ws.Cells.Item(37, 5).Value = '  -1.57%  '
$ws.Cells.Item(38, 5).Value = '  +2.38%  '
$ws.Cells.Item(39, 2).Value = 'ImmutableX'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(39, 4).Value = "'1.55"
$ws.Cells.Item(39, 5).Value = '  -3.17%  '
$ws.Cells.Item(40, 2).Value = 'Aptos'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(40, 4).Value = "'6.85"
$ws.Cells.Item(40, 5).Value = '  -1.17%  '
$ws.Cells.Item(41, 4).Value = "'164.10"
$ws.Cells.Item(41, 5).Value = '  -4.64%  '
$ws.Cells.Item(42, 4).Value = "'0.0783"
$ws.Cells.Item(42, 5).Value = '  -2.83%  '
$ws.Cells.Item(43, 4).Value = "'0.810"
$ws.Cells.Item(43, 5).Value = '  -0.52%  '
$ws.Cells.Item(44, 4).Value = "'25.57"
$ws.Cells.Item(44, 5).Value = '  -2.81%  '
$ws.Cells.Item(45, 5).Value = '  -0.08%  '
$ws.Cells.Item(46, 4).Value = "'41.87"
$ws.Cells.Item(46, 5).Value = '  -1.05%  '
$ws.Cells.Item(47, 5).Value = '  -0.38%  '
$ws.Cells.Item(48, 4).Value = "'4.40"
$ws.Cells.Item(48, 5).Value = '  +0.38%  '
$ws.Cells.Item(49, 5).Value = '  -1.02%  '
$ws.Cells.Item(50, 4).Value = '2.473.16'
$ws.Cells.Item(50, 5).Value = '  -0.27%  '
$ws.Cells.Item(51, 4).Value = "'0.918"
$ws.Cells.Item(51, 5).Value = '  +1.77%  '
